# Auto-generated edit script: apply scheduled-runner value updates to Hyperion_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 813.1667
$ws.Range("I31").Value = 813.1667
$ws.Range("K31").Value = 2439.5001
$ws.Range("M31").Value = -2209.5001

$ws.Range("H51").Value = 7217.5
$ws.Range("I51").Value = 7766.6665
$ws.Range("K51").Value = 7766.6665
$ws.Range("M51").Value = -7282.6665

$ws.Range("H53").Value = 5196.048
$ws.Range("J53").Value = 13327.75
$ws.Range("L53").Value = 13327.75
$ws.Range("N53").Value = -14601.75

$ws.Range("H62").Value = 5877.1113
$ws.Range("I62").Value = 3899.1667
$ws.Range("K62").Value = 3899.1667
$ws.Range("M62").Value = -3275.1667

$ws.Range("H65").Value = 5877.1113
$ws.Range("I65").Value = 3899.1667
$ws.Range("K65").Value = 19495.8335
$ws.Range("M65").Value = -16375.8335

$ws.Range("H70").Value = 5464.294
$ws.Range("J70").Value = 6239.3257
$ws.Range("L70").Value = 18717.9771
$ws.Range("N70").Value = -19257.9771

$ws.Range("H73").Value = 5464.294
$ws.Range("J73").Value = 6239.3257
$ws.Range("L73").Value = 18717.9771
$ws.Range("N73").Value = -20589.9771

$ws.Range("H80").Value = 592.88464
$ws.Range("J80").Value = 730.4286
$ws.Range("L80").Value = 2191.2858
$ws.Range("N80").Value = -4187.2858

$ws.Range("H83").Value = 592.88464
$ws.Range("J83").Value = 730.4286
$ws.Range("L83").Value = 6573.8574
$ws.Range("N83").Value = -16557.8574

$ws.Range("H98").Value = 1313.5555
$ws.Range("I98").Value = 1332
$ws.Range("K98").Value = 1332
$ws.Range("M98").Value = 166

$ws.Range("H106").Value = 3164.6667
$ws.Range("I106").Value = 3164.6667
$ws.Range("K106").Value = 3164.6667
$ws.Range("M106").Value = -2533.6667

$ws.Range("H118").Value = 100000430
$ws.Range("I118").Value = 100000430
$ws.Range("K118").Value = 300001290
$ws.Range("M118").Value = -299999633

$ws.Range("H122").Value = 1313.5555
$ws.Range("I122").Value = 1332
$ws.Range("K122").Value = 3996
$ws.Range("M122").Value = -1546

$ws.Range("H132").Value = 30306140
$ws.Range("I132").Value = 43481852
$ws.Range("J132").Value = 2007.5
$ws.Range("K132").Value = 130445556
$ws.Range("L132").Value = 6022.5
$ws.Range("M132").Value = -130443026
$ws.Range("N132").Value = -11082.5

$ws.Range("H137").Value = 59889.484
$ws.Range("I137").Value = 86384.48
$ws.Range("J137").Value = 4250
$ws.Range("K137").Value = 259153.44
$ws.Range("L137").Value = 12750
$ws.Range("M137").Value = -256603.44
$ws.Range("N137").Value = -17850

$ws.Range("H138").Value = 2769.7058
$ws.Range("J138").Value = 4832
$ws.Range("L138").Value = 14496
$ws.Range("N138").Value = -24776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 956.7273
$ws.Range("J4").Value = 538.5
$ws.Range("L4").Value = 538.5
$ws.Range("N4").Value = -770.5

$ws.Range("H32").Value = 2806.2207
$ws.Range("I32").Value = 1578.3135
$ws.Range("K32").Value = 1578.3135
$ws.Range("M32").Value = -1291.3135

$ws.Range("H45").Value = 14390612
$ws.Range("I45").Value = 23977938
$ws.Range("J45").Value = 9623.25
$ws.Range("K45").Value = 23977938
$ws.Range("L45").Value = 9623.25
$ws.Range("M45").Value = -23977561
$ws.Range("N45").Value = -10377.25

$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("M46").Value = -4681

$ws.Range("H74").Value = 74291.55499999999
$ws.Range("I74").Value = 5147.619
$ws.Range("K74").Value = 5147.619
$ws.Range("M74").Value = -4273.619

$ws.Range("H77").Value = 74291.55499999999
$ws.Range("I77").Value = 5147.619
$ws.Range("K77").Value = 25738.095
$ws.Range("M77").Value = -21370.095

$ws.Range("H122").Value = 803911.9
$ws.Range("I122").Value = 2380.7144
$ws.Range("K122").Value = 7142.1432
$ws.Range("M122").Value = -4692.1432

$ws.Range("H132").Value = 2297.457
$ws.Range("I132").Value = 2012.64
$ws.Range("J132").Value = 3009.5
$ws.Range("K132").Value = 6037.92
$ws.Range("L132").Value = 9028.5
$ws.Range("M132").Value = -3507.92
$ws.Range("N132").Value = -14088.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 31251000
$ws.Range("I105").Value = 31251000
$ws.Range("K105").Value = 31251000
$ws.Range("M105").Value = -31249253

$ws.Range("H134").Value = 3129.6316
$ws.Range("I134").Value = 1435.3438
$ws.Range("J134").Value = 12165.833
$ws.Range("K134").Value = 4306.0314
$ws.Range("L134").Value = 36497.499
$ws.Range("M134").Value = -1771.0314
$ws.Range("N134").Value = -41567.499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14838.432
$ws.Range("I31").Value = 1786.3704
$ws.Range("J31").Value = 21364.463
$ws.Range("K31").Value = 1786.3704
$ws.Range("L31").Value = 21364.463
$ws.Range("M31").Value = -1491.3704
$ws.Range("N31").Value = -21954.463

$ws.Range("H34").Value = 14838.432
$ws.Range("I34").Value = 1786.3704
$ws.Range("J34").Value = 21364.463
$ws.Range("K34").Value = 1786.3704
$ws.Range("L34").Value = 21364.463
$ws.Range("M34").Value = -1584.3704
$ws.Range("N34").Value = -21768.463

$ws.Range("H35").Value = 3579.8
$ws.Range("I35").Value = 1724.75
$ws.Range("J35").Value = 11000
$ws.Range("K35").Value = 1724.75
$ws.Range("L35").Value = 11000
$ws.Range("M35").Value = -1430.75
$ws.Range("N35").Value = -11588

$ws.Range("H58").Value = 2245.838
$ws.Range("I58").Value = 1910.8462
$ws.Range("J58").Value = 3037.6365
$ws.Range("K58").Value = 1910.8462
$ws.Range("L58").Value = 3037.6365
$ws.Range("M58").Value = -1707.8462
$ws.Range("N58").Value = -3443.6365

$ws.Range("H99").Value = 3677.1428
$ws.Range("I99").Value = 2935.625
$ws.Range("K99").Value = 2935.625
$ws.Range("M99").Value = -1437.625

$ws.Range("H122").Value = 3373.1667
$ws.Range("I122").Value = 2396.182
$ws.Range("J122").Value = 4199.846
$ws.Range("K122").Value = 7188.545999999999
$ws.Range("L122").Value = 12599.538
$ws.Range("M122").Value = -4738.545999999999
$ws.Range("N122").Value = -17499.538

$ws.Range("H126").Value = 3677.1428
$ws.Range("I126").Value = 2935.625
$ws.Range("K126").Value = 8806.875
$ws.Range("M126").Value = -6336.875

$ws.Range("H132").Value = 92459.19
$ws.Range("I132").Value = 69183.87
$ws.Range("K132").Value = 207551.61
$ws.Range("M132").Value = -205021.61

$ws.Range("H134").Value = 2701.2942
$ws.Range("I134").Value = 1411.0834
$ws.Range("J134").Value = 5797.8
$ws.Range("K134").Value = 4233.2502
$ws.Range("L134").Value = 17393.4
$ws.Range("M134").Value = -1698.2502
$ws.Range("N134").Value = -22463.4

$ws.Range("H136").Value = 2245.838
$ws.Range("I136").Value = 1910.8462
$ws.Range("J136").Value = 3037.6365
$ws.Range("K136").Value = 5732.5386
$ws.Range("L136").Value = 9112.9095
$ws.Range("M136").Value = -3182.5386
$ws.Range("N136").Value = -14212.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 479.33334
$ws.Range("I6").Value = 442.14285
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 1326.42855
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -1213.42855
$ws.Range("N6").Value = -3226

$ws.Range("H23").Value = 296.91666
$ws.Range("I23").Value = 157
$ws.Range("J23").Value = 343.55554
$ws.Range("K23").Value = 471
$ws.Range("L23").Value = 1030.66662
$ws.Range("M23").Value = -236
$ws.Range("N23").Value = -1500.66662

$ws.Range("H68").Value = 668.375
$ws.Range("I68").Value = 621
$ws.Range("K68").Value = 1863
$ws.Range("M68").Value = -1052

$ws.Range("H71").Value = 668.375
$ws.Range("I71").Value = 621
$ws.Range("K71").Value = 5589
$ws.Range("M71").Value = -1533

$ws.Range("H132").Value = 85123.586
$ws.Range("I132").Value = 126148
$ws.Range("K132").Value = 1135332
$ws.Range("M132").Value = -1132802

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 750.7857
$ws.Range("I2").Value = 1158.5555
$ws.Range("K2").Value = 1158.5555
$ws.Range("M2").Value = -1045.5555

$ws.Range("H80").Value = 1357842.6
$ws.Range("I80").Value = 3050527.5
$ws.Range("J80").Value = 3694.8
$ws.Range("K80").Value = 3050527.5
$ws.Range("L80").Value = 3694.8
$ws.Range("M80").Value = -3049529.5
$ws.Range("N80").Value = -5690.8

$ws.Range("H83").Value = 1357842.6
$ws.Range("I83").Value = 3050527.5
$ws.Range("J83").Value = 3694.8
$ws.Range("K83").Value = 15252637.5
$ws.Range("L83").Value = 18474
$ws.Range("M83").Value = -15247645.5
$ws.Range("N83").Value = -28458

$ws.Range("H132").Value = 3862.9395
$ws.Range("I132").Value = 3378.3845
$ws.Range("K132").Value = 10135.1535
$ws.Range("M132").Value = -7605.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 77354.81
$ws.Range("I136").Value = 120504.88
$ws.Range("J136").Value = 3999.7
$ws.Range("K136").Value = 361514.64
$ws.Range("L136").Value = 11999.1
$ws.Range("M136").Value = -358964.64
$ws.Range("N136").Value = -17099.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1676.9231
$ws.Range("I122").Value = 1300.2273
$ws.Range("K122").Value = 3900.6819
$ws.Range("M122").Value = -1450.6819

$ws.Range("H136").Value = 3165.5454
$ws.Range("I136").Value = 2427.25
$ws.Range("J136").Value = 7300
$ws.Range("K136").Value = 7281.75
$ws.Range("L136").Value = 21900
$ws.Range("M136").Value = -4731.75
$ws.Range("N136").Value = -27000
